$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -20.01078349860769
$ws.Range("C2").Value = 1.96280632347657
$ws.Range("D2").Value = -20.01078349860769
$ws.Range("E2").Value = -20.01078349860769
$ws.Range("F2").Value = -20.01078349860769
$ws.Range("G2").Value = -20.01078349860769
$ws.Range("H2").Value = -20.01078349860769
$ws.Range("I2").Value = -20.01078349860769
$ws.Range("J2").Value = -20.01078349860769
$ws.Range("K2").Value = -20.01078349860769
$ws.Range("B3").Value = -20.01078349860769
$ws.Range("C3").Value = -20.01078349860769
$ws.Range("D3").Value = -20.01078349860769
$ws.Range("E3").Value = -20.01078349860769
$ws.Range("F3").Value = -20.01078349860769
$ws.Range("G3").Value = -20.01078349860769
$ws.Range("H3").Value = -20.01078349860769
$ws.Range("I3").Value = 0.983114742351349
$ws.Range("J3").Value = -20.01078349860769
$ws.Range("K3").Value = -20.01078349860769
$ws.Range("B4").Value = -20.01078349860769
$ws.Range("C4").Value = 1.988086821132749
$ws.Range("D4").Value = 1.659890469468474
$ws.Range("E4").Value = -20.01078349860769
$ws.Range("F4").Value = 3.342536330011887
$ws.Range("G4").Value = -20.01078349860769
$ws.Range("H4").Value = 1.739946791627218
$ws.Range("I4").Value = -20.01078349860769
$ws.Range("J4").Value = 0.7248725070046795
$ws.Range("K4").Value = -20.01078349860769
$ws.Range("B5").Value = -20.01078349860769
$ws.Range("C5").Value = 1.66376611911916
$ws.Range("D5").Value = -20.01078349860769
$ws.Range("E5").Value = -20.01078349860769
$ws.Range("F5").Value = -20.01078349860769
$ws.Range("G5").Value = 4.321926797551959
$ws.Range("H5").Value = -20.01078349860769
$ws.Range("I5").Value = -20.01078349860769
$ws.Range("J5").Value = -20.01078349860769
$ws.Range("K5").Value = -20.01078349860769
$ws.Range("B6").Value = -20.01078349860769
$ws.Range("C6").Value = -20.01078349860769
$ws.Range("D6").Value = -20.01078349860769
$ws.Range("E6").Value = -20.01078349860769
$ws.Range("F6").Value = -20.01078349860769
$ws.Range("G6").Value = -20.01078349860769
$ws.Range("H6").Value = -20.01078349860769
$ws.Range("I6").Value = -20.01078349860769
$ws.Range("J6").Value = -20.01078349860769
$ws.Range("K6").Value = -20.01078349860769
$ws.Range("B7").Value = 2.439995071971311
$ws.Range("C7").Value = -20.01078349860769
$ws.Range("D7").Value = -20.01078349860769
$ws.Range("E7").Value = -20.01078349860769
$ws.Range("F7").Value = -20.01078349860769
$ws.Range("G7").Value = -20.01078349860769
$ws.Range("H7").Value = -20.01078349860769
$ws.Range("I7").Value = -20.01078349860769
$ws.Range("J7").Value = -20.01078349860769
$ws.Range("K7").Value = -20.01078349860769
$ws.Range("B8").Value = -20.01078349860769
$ws.Range("C8").Value = -20.01078349860769
$ws.Range("D8").Value = -20.01078349860769
$ws.Range("E8").Value = 1.840026363784017
$ws.Range("F8").Value = -20.01078349860769
$ws.Range("G8").Value = -20.01078349860769
$ws.Range("H8").Value = -20.01078349860769
$ws.Range("I8").Value = -20.01078349860769
$ws.Range("J8").Value = -20.01078349860769
$ws.Range("K8").Value = -20.01078349860769
$ws.Range("B9").Value = 3.865283828425842
$ws.Range("C9").Value = -20.01078349860769
$ws.Range("D9").Value = -20.01078349860769
$ws.Range("E9").Value = -20.01078349860769
$ws.Range("F9").Value = -20.01078349860769
$ws.Range("G9").Value = -20.01078349860769
$ws.Range("H9").Value = -20.01078349860769
$ws.Range("I9").Value = -20.01078349860769
$ws.Range("J9").Value = -20.01078349860769
$ws.Range("K9").Value = -20.01078349860769
$ws.Range("B10").Value = -20.01078349860769
$ws.Range("C10").Value = -20.01078349860769
$ws.Range("D10").Value = -20.01078349860769
$ws.Range("E10").Value = -20.01078349860769
$ws.Range("F10").Value = -20.01078349860769
$ws.Range("G10").Value = -20.01078349860769
$ws.Range("H10").Value = -20.01078349860769
$ws.Range("I10").Value = 1.839484618395143
$ws.Range("J10").Value = -20.01078349860769
$ws.Range("K10").Value = 2.303967064886773
$ws.Range("B11").Value = -20.01078349860769
$ws.Range("C11").Value = -20.01078349860769
$ws.Range("D11").Value = -20.01078349860769
$ws.Range("E11").Value = 2.889513786163938
$ws.Range("F11").Value = -20.01078349860769
$ws.Range("G11").Value = -20.01078349860769
$ws.Range("H11").Value = -20.01078349860769
$ws.Range("I11").Value = -20.01078349860769
$ws.Range("J11").Value = -20.01078349860769
$ws.Range("K11").Value = 1.788950680923222
$ws.Range("B12").Value = -20.01078349860769
$ws.Range("C12").Value = -20.01078349860769
$ws.Range("D12").Value = -20.01078349860769
$ws.Range("E12").Value = -20.01078349860769
$ws.Range("F12").Value = -20.01078349860769
$ws.Range("G12").Value = -20.01078349860769
$ws.Range("H12").Value = -20.01078349860769
$ws.Range("I12").Value = -20.01078349860769
$ws.Range("J12").Value = -20.01078349860769
$ws.Range("K12").Value = -20.01078349860769
$ws.Range("B13").Value = -20.01078349860769
$ws.Range("C13").Value = -20.01078349860769
$ws.Range("D13").Value = -20.01078349860769
$ws.Range("E13").Value = 2.527791237364887
$ws.Range("F13").Value = -20.01078349860769
$ws.Range("G13").Value = -20.01078349860769
$ws.Range("H13").Value = -20.01078349860769
$ws.Range("I13").Value = -20.01078349860769
$ws.Range("J13").Value = 1.538272107308043
$ws.Range("K13").Value = 1.606052925832458
$ws.Range("B14").Value = -20.01078349860769
$ws.Range("C14").Value = -20.01078349860769
$ws.Range("D14").Value = 1.544427693558402
$ws.Range("E14").Value = -20.01078349860769
$ws.Range("F14").Value = -20.01078349860769
$ws.Range("G14").Value = -20.01078349860769
$ws.Range("H14").Value = -20.01078349860769
$ws.Range("I14").Value = -20.01078349860769
$ws.Range("J14").Value = -20.01078349860769
$ws.Range("K14").Value = 1.880147579176696
$ws.Range("B15").Value = -20.01078349860769
$ws.Range("C15").Value = -20.01078349860769
$ws.Range("D15").Value = 1.757057047117734
$ws.Range("E15").Value = -20.01078349860769
$ws.Range("F15").Value = -20.01078349860769
$ws.Range("G15").Value = -20.01078349860769
$ws.Range("H15").Value = -20.01078349860769
$ws.Range("I15").Value = -20.01078349860769
$ws.Range("J15").Value = -20.01078349860769
$ws.Range("K15").Value = -20.01078349860769
$ws.Range("B16").Value = -20.01078349860769
$ws.Range("C16").Value = -20.01078349860769
$ws.Range("D16").Value = -20.01078349860769
$ws.Range("E16").Value = -20.01078349860769
$ws.Range("F16").Value = -20.01078349860769
$ws.Range("G16").Value = -20.01078349860769
$ws.Range("H16").Value = -20.01078349860769
$ws.Range("I16").Value = -20.01078349860769
$ws.Range("J16").Value = 1.998390749043498
$ws.Range("K16").Value = -20.01078349860769
$ws.Range("B17").Value = -20.01078349860769
$ws.Range("C17").Value = 2.11712803716661
$ws.Range("D17").Value = 1.827129647208195
$ws.Range("E17").Value = -20.01078349860769
$ws.Range("F17").Value = -20.01078349860769
$ws.Range("G17").Value = -20.01078349860769
$ws.Range("H17").Value = 2.049296575418105
$ws.Range("I17").Value = 2.097184763763146
$ws.Range("J17").Value = 2.576498901663761
$ws.Range("K17").Value = -20.01078349860769
$ws.Range("B18").Value = -20.01078349860769
$ws.Range("C18").Value = -20.01078349860769
$ws.Range("D18").Value = -20.01078349860769
$ws.Range("E18").Value = -20.01078349860769
$ws.Range("F18").Value = -20.01078349860769
$ws.Range("G18").Value = -20.01078349860769
$ws.Range("H18").Value = 1.818420876104589
$ws.Range("I18").Value = 2.067134323029873
$ws.Range("J18").Value = 2.454787140812163
$ws.Range("K18").Value = -20.01078349860769
$ws.Range("B19").Value = -20.01078349860769
$ws.Range("C19").Value = -20.01078349860769
$ws.Range("D19").Value = 2.038859352716904
$ws.Range("E19").Value = -20.01078349860769
$ws.Range("F19").Value = -20.01078349860769
$ws.Range("G19").Value = -20.01078349860769
$ws.Range("H19").Value = 1.689372129142714
$ws.Range("I19").Value = 1.81905573546812
$ws.Range("J19").Value = -20.01078349860769
$ws.Range("K19").Value = -20.01078349860769
$ws.Range("B20").Value = -20.01078349860769
$ws.Range("C20").Value = 1.082260679889003
$ws.Range("D20").Value = 1.528866159939828
$ws.Range("E20").Value = -20.01078349860769
$ws.Range("F20").Value = 3.301018715375643
$ws.Range("G20").Value = -20.01078349860769
$ws.Range("H20").Value = 1.657000701375037
$ws.Range("I20").Value = 1.290815958585341
$ws.Range("J20").Value = -20.01078349860769
$ws.Range("K20").Value = 2.287149533449363
$ws.Range("B21").Value = -20.01078349860769
$ws.Range("C21").Value = 1.32804046045171
$ws.Range("D21").Value = -20.01078349860769
$ws.Range("E21").Value = 1.697209575278437
$ws.Range("F21").Value = -20.01078349860769
$ws.Range("G21").Value = -20.01078349860769
$ws.Range("H21").Value = 1.386812280753763
$ws.Range("I21").Value = -20.01078349860769
$ws.Range("J21").Value = -20.01078349860769
$ws.Range("K21").Value = -20.01078349860769
